$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: fix typo in options string "[_1,_2]" -> "[1,2]" ---
$ws.Range("B6").Value = "[1,2]"

# --- Row 9: was a blank "quote-prefixed" question cell; give it real text + an answer ---
$ws.Range("A9").Value = "Please tell us how interested you are in the NHL."
$ws.Range("B9").Value = 1

# --- Row 10: replace the Hispanic/Latino question text with the new race & ethnicity
#     consent intro, make its answer cell bold-free but theme-colored font, taller row ---
$ws.Range("A10").Value = "The next question will be about race and ethnicity. A " + [char]8220 + "Prefer not to answer" + [char]8221 + " option is available for you to select, at your discretion. Collecting such information enables us to provide a more refined research analysis.
Participation is always voluntary, and your responses are used for research purposes only, combined with the answers from all other participants. We will provide our client only anonymous, aggregated results. The data will be held for no longer than 12 months.
Do you accept the collection of race and ethnicity related data?
Select only one"
$ws.Rows(10).RowHeight = 87

$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = 1
$ws.Range("B10").Font.ThemeColor = 1

# --- New row 11: subscription streaming question ---
$ws.Range("A11").Value = "Which of the following subscription streaming services do you subscribe to?"
$ws.Range("A1").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B11").Value = 1
$ws.Rows(11).RowHeight = 19.5

# --- New row 12: original Hispanic/Latino question text, now prefixed with a
#     sensitive-topic disclaimer, moved down here ---
$ws.Range("A12").Value = "This is a topic of a sensitive nature. Answering is voluntary, however, collecting such information enables us to provide a more refined research analysis.
Are you of Hispanic, Latino or Spanish origin?
If you don" + [char]8217 + "t agree to provide us such information, a " + [char]8220 + "Prefer not to answer" + [char]8221 + " option is available for you to select, at your discretion.
For any survey research purposes, your responses are combined with the answers from all other participants. We will provide our client only anonymous results, unless you separately consent otherwise. The data will be held by us for the research purposes no longer than 12 months."
$ws.Range("A1").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Value = 1
$ws.Rows(12).RowHeight = 60
